$d = $word.ActiveDocument

# Replace "Summary 2" header with "Summary 0"
$d.Content.Find.Execute("Summary 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Summary 0", 2)

# Replace "N 1" header with "Missing 1"
$d.Content.Find.Execute("N 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Missing 1", 2)

# Replace "N 2" header with "Missing 0"
$d.Content.Find.Execute("N 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Missing 0", 2)
